$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4025
$ws.Range("I74").Value = 3685.7144
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 3685.7144
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -2749.7144
$ws.Range("N74").Value = -6372
$ws.Range("H77").Value = 4025
$ws.Range("I77").Value = 3685.7144
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 18428.572
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -13748.572
$ws.Range("N77").Value = -31860
$ws.Range("H100").Value = 11112420
$ws.Range("I100").Value = 16668068
$ws.Range("J100").Value = 1126
$ws.Range("K100").Value = 16668068
$ws.Range("L100").Value = 1126
$ws.Range("M100").Value = -16667527
$ws.Range("N100").Value = -2208
$ws.Range("H112").Value = 9741197
$ws.Range("J112").Value = 10910040
$ws.Range("L112").Value = 32730120
$ws.Range("N112").Value = -32732336
$ws.Range("H129").Value = 1116.6765
$ws.Range("I129").Value = 379.1111
$ws.Range("J129").Value = 1382.2
$ws.Range("K129").Value = 1137.3333
$ws.Range("L129").Value = 4146.6
$ws.Range("M129").Value = 3862.6667
$ws.Range("N129").Value = -14146.6
$ws.Range("H133").Value = 23119.875
$ws.Range("J133").Value = 23119.875
$ws.Range("L133").Value = 23119.875
$ws.Range("N133").Value = -33239.875
$ws.Range("H137").Value = 71430300
$ws.Range("I137").Value = 142858050
$ws.Range("J137").Value = 2545
$ws.Range("K137").Value = 428574150
$ws.Range("L137").Value = 7635
$ws.Range("M137").Value = -428571600
$ws.Range("N137").Value = -12735
$ws.Range("H138").Value = 2387.7937
$ws.Range("I138").Value = 864.10345
$ws.Range("J138").Value = 3037.603
$ws.Range("K138").Value = 2592.31035
$ws.Range("L138").Value = 9112.809000000001
$ws.Range("M138").Value = 2547.68965
$ws.Range("N138").Value = -19392.809

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 43170.043
$ws.Range("I2").Value = 60305.47
$ws.Range("J2").Value = 1555.4286
$ws.Range("K2").Value = 60305.47
$ws.Range("L2").Value = 1555.4286
$ws.Range("M2").Value = -60192.47
$ws.Range("N2").Value = -1781.4286
$ws.Range("H32").Value = 17828.734
$ws.Range("I32").Value = 2169.758
$ws.Range("K32").Value = 2169.758
$ws.Range("M32").Value = -1882.758
$ws.Range("H63").Value = 14450
$ws.Range("I63").Value = 14450
$ws.Range("K63").Value = 14450
$ws.Range("M63").Value = -13764
$ws.Range("H66").Value = 14450
$ws.Range("I66").Value = 14450
$ws.Range("K66").Value = 72250
$ws.Range("M66").Value = -68818
$ws.Range("H74").Value = 3378.3455
$ws.Range("I74").Value = 1008.2195
$ws.Range("J74").Value = 10319.429
$ws.Range("K74").Value = 1008.2195
$ws.Range("L74").Value = 10319.429
$ws.Range("M74").Value = -134.2195
$ws.Range("N74").Value = -12067.429
$ws.Range("H77").Value = 3378.3455
$ws.Range("I77").Value = 1008.2195
$ws.Range("J77").Value = 10319.429
$ws.Range("K77").Value = 5041.0975
$ws.Range("L77").Value = 51597.145
$ws.Range("M77").Value = -673.0974999999999
$ws.Range("N77").Value = -60333.145
$ws.Range("H110").Value = 1571
$ws.Range("I110").Value = 1050
$ws.Range("J110").Value = 2613
$ws.Range("K110").Value = 1050
$ws.Range("L110").Value = 2613
$ws.Range("M110").Value = 995
$ws.Range("N110").Value = -6703
$ws.Range("H116").Value = 43170.043
$ws.Range("I116").Value = 60305.47
$ws.Range("J116").Value = 1555.4286
$ws.Range("K116").Value = 60305.47
$ws.Range("L116").Value = 1555.4286
$ws.Range("M116").Value = -58011.47
$ws.Range("N116").Value = -6143.4286
$ws.Range("H122").Value = 2300.348
$ws.Range("I122").Value = 2110
$ws.Range("J122").Value = 2985.6
$ws.Range("K122").Value = 6330
$ws.Range("L122").Value = 8956.799999999999
$ws.Range("M122").Value = -3880
$ws.Range("N122").Value = -13856.8
$ws.Range("H132").Value = 1595.4333
$ws.Range("I132").Value = 1229.9246
$ws.Range("J132").Value = 4362.857
$ws.Range("K132").Value = 3689.7738
$ws.Range("L132").Value = 13088.571
$ws.Range("M132").Value = -1159.7738
$ws.Range("N132").Value = -18148.571

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 43170.043
$ws.Range("I3").Value = 60305.47
$ws.Range("J3").Value = 1555.4286
$ws.Range("K3").Value = 60305.47
$ws.Range("L3").Value = 1555.4286
$ws.Range("M3").Value = -60191.47
$ws.Range("N3").Value = -1783.4286
$ws.Range("H105").Value = 3058.9167
$ws.Range("I105").Value = 2837.4583
$ws.Range("K105").Value = 2837.4583
$ws.Range("M105").Value = -1090.4583
$ws.Range("H107").Value = 743.82355
$ws.Range("I107").Value = 617.63635
$ws.Range("J107").Value = 975.1667
$ws.Range("K107").Value = 617.63635
$ws.Range("L107").Value = 975.1667
$ws.Range("M107").Value = 1302.36365
$ws.Range("N107").Value = -4815.1667
$ws.Range("H134").Value = 2990.7073
$ws.Range("I134").Value = 1957.0741
$ws.Range("J134").Value = 4984.143
$ws.Range("K134").Value = 5871.2223
$ws.Range("L134").Value = 14952.429
$ws.Range("M134").Value = -3336.2223
$ws.Range("N134").Value = -20022.429

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1787.3226
$ws.Range("I31").Value = 1046.3
$ws.Range("J31").Value = 3134.6365
$ws.Range("K31").Value = 1046.3
$ws.Range("L31").Value = 3134.6365
$ws.Range("M31").Value = -751.3
$ws.Range("N31").Value = -3724.6365
$ws.Range("H34").Value = 1787.3226
$ws.Range("I34").Value = 1046.3
$ws.Range("J34").Value = 3134.6365
$ws.Range("K34").Value = 1046.3
$ws.Range("L34").Value = 3134.6365
$ws.Range("M34").Value = -844.3
$ws.Range("N34").Value = -3538.6365
$ws.Range("H58").Value = 2139.7026
$ws.Range("I58").Value = 909.26086
$ws.Range("J58").Value = 4161.143
$ws.Range("K58").Value = 909.26086
$ws.Range("L58").Value = 4161.143
$ws.Range("M58").Value = -706.26086
$ws.Range("N58").Value = -4567.143
$ws.Range("H134").Value = 1987.0588
$ws.Range("I134").Value = 1117.1628
$ws.Range("J134").Value = 6662.75
$ws.Range("K134").Value = 3351.4884
$ws.Range("L134").Value = 19988.25
$ws.Range("M134").Value = -816.4884000000002
$ws.Range("N134").Value = -25058.25
$ws.Range("H136").Value = 2139.7026
$ws.Range("I136").Value = 909.26086
$ws.Range("J136").Value = 4161.143
$ws.Range("K136").Value = 2727.78258
$ws.Range("L136").Value = 12483.429
$ws.Range("M136").Value = -177.7825800000001
$ws.Range("N136").Value = -17583.429

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 298.35294
$ws.Range("I11").Value = 310.42856
$ws.Range("J11").Value = 242
$ws.Range("K11").Value = 931.28568
$ws.Range("L11").Value = 726
$ws.Range("M11").Value = -791.28568
$ws.Range("N11").Value = -1006
$ws.Range("H92").Value = 702.8
$ws.Range("I92").Value = 567.3333
$ws.Range("J92").Value = 760.8570999999999
$ws.Range("K92").Value = 1701.9999
$ws.Range("L92").Value = 2282.5713
$ws.Range("M92").Value = -453.9999
$ws.Range("N92").Value = -4778.5713

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("H132").Value = 2590.3
$ws.Range("I132").Value = 2475.7693
$ws.Range("J132").Value = 2996.3635
$ws.Range("K132").Value = 7427.3079
$ws.Range("L132").Value = 8989.0905
$ws.Range("M132").Value = -4897.3079
$ws.Range("N132").Value = -14049.0905
$ws.Range("N64").Value = $null
$ws.Range("N67").Value = $null

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2017.6666
$ws.Range("I93").Value = 1502
$ws.Range("K93").Value = 1502
$ws.Range("M93").Value = -254
$ws.Range("H122").Value = 3614.75
$ws.Range("I122").Value = 2066.6667
$ws.Range("J122").Value = 3887.9412
$ws.Range("K122").Value = 6200.000100000001
$ws.Range("L122").Value = 11663.8236
$ws.Range("M122").Value = -3750.000100000001
$ws.Range("N122").Value = -16563.8236
$ws.Range("H132").Value = 3489.1365
$ws.Range("I132").Value = 2518.8276
$ws.Range("J132").Value = 5365.067
$ws.Range("K132").Value = 7556.4828
$ws.Range("L132").Value = 16095.201
$ws.Range("M132").Value = -5026.4828
$ws.Range("N132").Value = -21155.201
$ws.Range("H134").Value = 16714.5
$ws.Range("J134").Value = 26429
$ws.Range("L134").Value = 26429
$ws.Range("N134").Value = -36569
$ws.Range("H136").Value = 4412.3687
$ws.Range("I136").Value = 2972.6
$ws.Range("J136").Value = 9811.5
$ws.Range("K136").Value = 8917.799999999999
$ws.Range("L136").Value = 29434.5
$ws.Range("M136").Value = -6367.799999999999
$ws.Range("N136").Value = -34534.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1580
$ws.Range("J8").Value = 1475
$ws.Range("L8").Value = 1475
$ws.Range("N8").Value = -1755
$ws.Range("H9").Value = 18666.334
$ws.Range("I9").Value = 50000
$ws.Range("J9").Value = 2999.5
$ws.Range("K9").Value = 50000
$ws.Range("L9").Value = 2999.5
$ws.Range("M9").Value = -49860
$ws.Range("N9").Value = -3279.5
$ws.Range("H10").Value = 2880
$ws.Range("J10").Value = 2880
$ws.Range("L10").Value = 2880
$ws.Range("N10").Value = -3218
$ws.Range("H92").Value = 34999.668
$ws.Range("J92").Value = 34999.668
$ws.Range("L92").Value = 34999.668
$ws.Range("N92").Value = -39991.668
$ws.Range("H122").Value = 253001
$ws.Range("I122").Value = 253001
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 759003
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -756553
$ws.Range("H126").Value = 78091.16
$ws.Range("I126").Value = 143343.58
$ws.Range("J126").Value = 1963.3334
$ws.Range("K126").Value = 430030.74
$ws.Range("L126").Value = 5890.0002
$ws.Range("M126").Value = -427560.74
$ws.Range("N126").Value = -10830.0002
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H136").Value = 6430224
$ws.Range("I136").Value = 7114190.5
$ws.Range("K136").Value = 21342571.5
$ws.Range("M136").Value = -21340021.5
$ws.Range("N122").Value = $null
$ws.Range("N127").Value = $null
